$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column retains its text formatting (values like "1.000" or "16.00"
# must not be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.715.03'
$ws.Range('E2').Value = '  +6.81%  '
$ws.Range('D3').Value = '1.809.08'
$ws.Range('E3').Value = '  +4.51%  '
$ws.Range('D5').Value = '251.20'
$ws.Range('E5').Value = '  +3.70%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '0.4959'
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('D8').Value = '0.2790'
$ws.Range('E8').Value = '  +7.20%  '
$ws.Range('D9').Value = '0.06382'
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('D10').Value = '1.804.49'
$ws.Range('E10').Value = '  +4.18%  '
$ws.Range('D11').Value = '16.79'
$ws.Range('E11').Value = '  +4.39%  '
$ws.Range('D12').Value = '0.07114'
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('D13').Value = '0.6459'
$ws.Range('E13').Value = '  +5.71%  '
$ws.Range('D14').Value = '4.702'
$ws.Range('E14').Value = '  +4.28%  '
$ws.Range('D15').Value = '81.92'
$ws.Range('E15').Value = '  +5.84%  '
$ws.Range('D16').Value = '28.692.92'
$ws.Range('E16').Value = '  +7.67%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '0.000007347'
$ws.Range('E18').Value = '  +2.14%  '
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').Value = '12.25'
$ws.Range('E20').Value = '  +6.81%  '
$ws.Range('D21').Value = '2.036.53'
$ws.Range('E21').Value = '  +3.96%  '
$ws.Range('D22').Value = '4.613'
$ws.Range('E22').Value = '  +3.80%  '
$ws.Range('D23').Value = '8.884'
$ws.Range('E23').Value = '  +3.65%  '
$ws.Range('D24').Value = '5.312'
$ws.Range('E24').Value = '  +3.53%  '
$ws.Range('D25').Value = '142.76'
$ws.Range('E25').Value = '  +2.84%  '
$ws.Range('D26').Value = '16.00'
$ws.Range('E26').Value = '  +4.38%  '
$ws.Range('D27').Value = '1.879'
$ws.Range('E27').Value = '  +4.79%  '
$ws.Range('D28').Value = '111.27'
$ws.Range('E28').Value = '  +4.80%  '
$ws.Range('D29').Value = '1.386'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('D30').Value = '4.181'
$ws.Range('E30').Value = '  +5.91%  '
$ws.Range('D31').Value = '0.08359'
$ws.Range('E31').Value = '  +4.50%  '
$ws.Range('D32').Value = '3.840'
$ws.Range('E32').Value = '  +4.32%  '
$ws.Range('D33').Value = '0.04958'
$ws.Range('E33').Value = '  +9.41%  '
$ws.Range('D34').Value = '1.092'
$ws.Range('E34').Value = '  +7.96%  '
$ws.Range('D35').Value = '0.6713'
$ws.Range('E35').Value = '  +7.29%  '
$ws.Range('D36').Value = '2.664'
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('D37').Value = '0.9595'
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('D38').Value = '2.639'
$ws.Range('E38').Value = '  +8.03%  '
$ws.Range('D39').Value = '2.146'
$ws.Range('E39').Value = '  +4.10%  '
$ws.Range('E40').Value = '  +6.04%  '
$ws.Range('D41').Value = '5.923'
$ws.Range('E41').Value = '  +4.85%  '
$ws.Range('D42').Value = '0.9999'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '101.05'
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').Value = '0.4114'
$ws.Range('E44').Value = '  +6.09%  '
$ws.Range('D45').Value = '7.236'
$ws.Range('E45').Value = '  +4.03%  '
$ws.Range('D46').Value = '0.1225'
$ws.Range('E46').Value = '  +5.46%  '
$ws.Range('D47').Value = '0.05495'
$ws.Range('E47').Value = '  +1.98%  '
$ws.Range('D48').Value = '8.156'
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('D49').Value = '31.34'
$ws.Range('E49').Value = '  +3.38%  '
$ws.Range('D50').Value = '1.302'
$ws.Range('E50').Value = '  +4.59%  '
$ws.Range('D51').Value = '0.3603'
$ws.Range('E51').Value = '  +6.22%  '
